$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$part1 = "Captain Reilly Jr., son of Captain Reilly Sr., followed his father" + [char]0x2019 + "s pirate career. His whole life revolved around escaping poverty after his father plundered the "
$part2 = "wrong ship. "

$r1 = $tr.InsertAfter($part1)
$r1.Text = $r1.Text

$r2 = $tr.InsertAfter($part2)
$r2.Text = $r2.Text
